$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so values like "1.000"
# or "26.507.96" are not auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.507.96"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "1.734.54"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "247.03"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "0.4887"
$ws.Range("E7").Value = "  +1.59%  "

$ws.Range("D8").Value = "0.2665"
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").Value = "0.06218"
$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").Value = "1.729.75"
$ws.Range("E10").Value = "  -0.60%  "

$ws.Range("D11").Value = "0.07028"
$ws.Range("E11").Value = "  -1.41%  "

$ws.Range("D12").Value = "15.66"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("D13").Value = "4.591"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").Value = "0.6083"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").Value = "77.45"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "26.513.60"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "0.000007327"
$ws.Range("E18").Value = "  +6.33%  "

$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").Value = "11.51"
$ws.Range("E20").Value = "  -2.28%  "

$ws.Range("D21").Value = "1.953.77"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").Value = "4.554"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").Value = "8.749"
$ws.Range("E23").Value = "  -1.78%  "

$ws.Range("D24").Value = "5.229"
$ws.Range("E24").Value = "  -2.34%  "

$ws.Range("D25").Value = "140.74"
$ws.Range("E25").Value = "  +3.58%  "

$ws.Range("D26").Value = "15.42"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "1.410"
$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").Value = "1.772"
$ws.Range("E28").Value = "  -2.17%  "

$ws.Range("D29").Value = "107.82"
$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("D30").Value = "4.016"
$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").Value = "0.08057"
$ws.Range("E31").Value = "  +2.11%  "

$ws.Range("D32").Value = "3.688"
$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("D33").Value = "0.04561"
$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("D34").Value = "1.0000"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "2.613"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").Value = "1.007"
$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("D37").Value = "0.6345"
$ws.Range("E37").Value = "  -0.56%  "

$ws.Range("D38").Value = "0.8973"
$ws.Range("E38").Value = "  -3.72%  "

$ws.Range("D39").Value = "2.022"
$ws.Range("E39").Value = "  +1.09%  "

$ws.Range("D40").Value = "2.399"

$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").Value = "0.01502"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").Value = "101.56"
$ws.Range("E43").Value = "  -9.44%  "

$ws.Range("D44").Value = "5.417"
$ws.Range("E44").Value = "  -6.01%  "

$ws.Range("D45").Value = "0.3884"
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").Value = "6.924"
$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("D48").Value = "0.05396"
$ws.Range("E48").Value = "  +1.19%  "

$ws.Range("D49").Value = "7.816"
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("D50").Value = "30.49"
$ws.Range("E50").Value = "  -1.14%  "

$ws.Range("D51").Value = "1.255"
$ws.Range("E51").Value = "  -0.16%  "
